# Lexus dataframe cleanup: remove duplicated rows (rows 142-145 were
# duplicates of earlier ISF/LX/LFA entries) and re-pack the remaining
# data for rows 109-141 so that the "A" index column and the
# year/model/security columns line up with the de-duplicated dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Excel row, A (index), B (year), C (make), D (model), E (security)
$rowData = @(
    @(109, 108, 2010, "Lexus", "LX",  "Smart Key"),
    @(110, 109, 2010, "Lexus", "GS",  "Std"),
    @(111, 110, 2010, "Lexus", "GX",  "Std"),
    @(112, 111, 2010, "Lexus", "RX",  "Std"),
    @(113, 112, 2010, "Lexus", "SC",  "Std"),
    @(114, 113, 2011, "Lexus", "ES",  "Smart Key"),
    @(115, 114, 2011, "Lexus", "GX",  "Smart Key"),
    @(116, 115, 2011, "Lexus", "LS",  "Smart Key"),
    @(117, 116, 2011, "Lexus", "RX",  "Smart Key"),
    @(118, 117, 2011, "Lexus", "GS",  "Smart Key"),
    @(119, 118, 2011, "Lexus", "IS",  "Smart Key"),
    @(120, 119, 2011, "Lexus", "ISF", "Smart Key"),
    @(121, 121, 2011, "Lexus", "LX",  "Smart Key"),
    @(122, 122, 2011, "Lexus", "SC",  "Std"),
    @(123, 123, 2012, "Lexus", "ES",  "Smart Key"),
    @(124, 124, 2012, "Lexus", "GX",  "Smart Key"),
    @(125, 125, 2012, "Lexus", "LS",  "Smart Key"),
    @(126, 126, 2012, "Lexus", "RX",  "Smart Key"),
    @(127, 127, 2012, "Lexus", "GS",  "Smart Key"),
    @(128, 128, 2012, "Lexus", "IS",  "Smart Key"),
    @(129, 129, 2012, "Lexus", "ISF", "Smart Key"),
    @(130, 131, 2012, "Lexus", "LX",  "Smart Key"),
    @(131, 132, 2012, "Lexus", "LFA", "Std"),
    @(132, 133, 2012, "Lexus", "SC",  "Std"),
    @(133, 134, 2013, "Lexus", "ES",  "Smart Key"),
    @(134, 135, 2013, "Lexus", "GX",  "Smart Key"),
    @(135, 136, 2013, "Lexus", "LS",  "Smart Key"),
    @(136, 137, 2013, "Lexus", "RX",  "Smart Key"),
    @(137, 138, 2013, "Lexus", "GS",  "Smart Key"),
    @(138, 139, 2013, "Lexus", "IS",  "Smart Key"),
    @(139, 140, 2013, "Lexus", "ISF", "Smart Key"),
    @(140, 142, 2013, "Lexus", "LX",  "Smart Key"),
    @(141, 143, 2013, "Lexus", "LFA", "Std")
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
    $ws.Cells.Item($r, 5).Value = $entry[5]
}

# The last four rows (142-145) were exact duplicates of earlier rows and
# are removed entirely, shrinking the used range to A1:E141.
$ws.Rows("142:145").Delete()
